# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 56 (pushing the existing rows
# 56-74 down to 57-75), matching the new dimension A1:T75.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 56..74 down by one to make room for the new weekly record.
$ws.Rows(56).Insert()

# Populate the newly inserted row 56 with this week's data.
$ws.Range("A56").Value = 10
$ws.Range("B56").Value = "Vega Modelo de Temuco"
$ws.Range("C56").Value = "La Araucanía"
$ws.Range("D56").Value = 44992
$ws.Range("E56").Value = 9
$ws.Range("F56").Value = "Fruta"
$ws.Range("G56").Value = 100107
$ws.Range("H56").Value = "Otros"
$ws.Range("I56").Value = 100107011
$ws.Range("J56").Value = "Tuna"
$ws.Range("K56").Value = "Sin especificar"
$ws.Range("L56").Value = "Primera"
$ws.Range("M56").Value = 80
$ws.Range("N56").Value = 20000
$ws.Range("O56").Value = 20000
$ws.Range("P56").Value = 20000
$ws.Range("Q56").Value = '$/caja 16 kilos'
$ws.Range("R56").Value = "Provincia de Los Andes"
$ws.Range("S56").Value = 1250
$ws.Range("T56").Value = 16
